$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the date-formatting style used by the existing "Fecha" (D) column
$dateFormat = $ws.Range("D2").NumberFormat()

# Updated / new data rows 19-38 (rows shifted down, one new row inserted at 22,
# four brand-new rows appended at 35-38)
$rows = @(
    @{R=19; Y=2021; Mo=10; Da=14; H='Sin especificar'; I='Banquete'; J=630; K=1500; L=1600; M=1556; N='$/kilo'; O='Provincia de Linares'; P=1556; Q=1}
    @{R=20; Y=2021; Mo=10; Da=14; H='Sin especificar'; I='Primera'; J=480; K=1300; L=1400; M=1352; N='$/kilo'; O='Provincia de Linares'; P=1352; Q=1}
    @{R=21; Y=2021; Mo=10; Da=14; H='Sin especificar'; I='Segunda'; J=250; K=1100; L=1200; M=1160; N='$/kilo'; O='Provincia de Linares'; P=1160; Q=1}
    @{R=22; Y=2021; Mo=10; Da=14; H='Sin especificar'; I='Tercera'; J=300; K=800; L=900; M=867; N='$/kilo'; O='Provincia de Linares'; P=867; Q=1}
    @{R=23; Y=2021; Mo=9; Da=28; H='Sin especificar'; I='Banquete'; J=100; K=2000; L=2000; M=2000; N='$/caja 10 kilos'; O='Provincia de Linares'; P=200; Q=10}
    @{R=24; Y=2021; Mo=9; Da=28; H='Sin especificar'; I='Primera'; J=150; K=1500; L=1500; M=1500; N='$/caja 10 kilos'; O='Provincia de Linares'; P=150; Q=10}
    @{R=25; Y=2021; Mo=9; Da=28; H='Sin especificar'; I='Segunda'; J=50; K=1300; L=1300; M=1300; N='$/caja 10 kilos'; O='Provincia de Linares'; P=130; Q=10}
    @{R=26; Y=2020; Mo=11; Da=26; H='Verde'; I='Primera'; J=4300; K=1000; L=1000; M=1000; N='$/kilo'; O='Provincia de Linares'; P=1000; Q=1}
    @{R=27; Y=2020; Mo=11; Da=26; H='Verde'; I='Segunda'; J=2500; K=800; L=800; M=800; N='$/kilo'; O='Provincia de Linares'; P=800; Q=1}
    @{R=28; Y=2021; Mo=10; Da=5; H='Sin especificar'; I='Banquete'; J=780; K=1500; L=1600; M=1558; N='$/kilo'; O='Provincia de Linares'; P=1558; Q=1}
    @{R=29; Y=2021; Mo=10; Da=5; H='Sin especificar'; I='Primera'; J=520; K=1300; L=1400; M=1348; N='$/kilo'; O='Provincia de Linares'; P=1348; Q=1}
    @{R=30; Y=2021; Mo=10; Da=5; H='Sin especificar'; I='Segunda'; J=400; K=1000; L=1200; M=1100; N='$/kilo'; O='Provincia de Linares'; P=1100; Q=1}
    @{R=31; Y=2020; Mo=11; Da=24; H='Verde'; I='Primera'; J=4300; K=1000; L=1000; M=1000; N='$/kilo'; O='Región Metropolitana'; P=1000; Q=1}
    @{R=32; Y=2020; Mo=11; Da=24; H='Verde'; I='Segunda'; J=2500; K=800; L=800; M=800; N='$/kilo'; O='Región Metropolitana'; P=800; Q=1}
    @{R=33; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Banquete'; J=1300; K=1500; L=1600; M=1554; N='$/kilo'; O='Provincia de Linares'; P=1554; Q=1}
    @{R=34; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Banquete'; J=700; K=1400; L=1500; M=1457; N='$/kilo'; O='Región Metropolitana'; P=1457; Q=1}
    @{R=35; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Primera'; J=900; K=1300; L=1400; M=1356; N='$/kilo'; O='Provincia de Linares'; P=1356; Q=1}
    @{R=36; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Primera'; J=500; K=1200; L=1300; M=1260; N='$/kilo'; O='Región Metropolitana'; P=1260; Q=1}
    @{R=37; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Segunda'; J=500; K=1100; L=1200; M=1160; N='$/kilo'; O='Provincia de Linares'; P=1160; Q=1}
    @{R=38; Y=2021; Mo=10; Da=7; H='Sin especificar'; I='Segunda'; J=200; K=1000; L=1100; M=1050; N='$/kilo'; O='Región Metropolitana'; P=1050; Q=1}
)

foreach ($row in $rows) {
    $r = $row.R
    $addr = "D" + $r
    if ($r -gt 34) {
        # Brand-new rows have no pre-existing style, so carry over the date format explicitly
        $ws.Range($addr).NumberFormat = $dateFormat
    }
    $ws.Range($addr).Value = (Get-Date -Year $row.Y -Month $row.Mo -Day $row.Da -Hour 0 -Minute 0 -Second 0)
    $ws.Range("H" + $r).Value = $row.H
    $ws.Range("I" + $r).Value = $row.I
    $ws.Range("J" + $r).Value = $row.J
    $ws.Range("K" + $r).Value = $row.K
    $ws.Range("L" + $r).Value = $row.L
    $ws.Range("M" + $r).Value = $row.M
    $ws.Range("N" + $r).Value = $row.N
    $ws.Range("O" + $r).Value = $row.O
    $ws.Range("P" + $r).Value = $row.P
    $ws.Range("Q" + $r).Value = $row.Q

    if ($r -gt 34) {
        # Columns A, B, C, E, F, G and R are constant across every record in this sheet
        $ws.Range("A" + $r).Value = 6
        $ws.Range("B" + $r).Value = 'Mercado Mayorista Lo Valledor de Santiago'
        $ws.Range("C" + $r).Value = 'Metropolitana'
        $ws.Range("E" + $r).Value = 13
        $ws.Range("F" + $r).Value = 300000000
        $ws.Range("G" + $r).Value = 'Espárragos'
        $ws.Range("R" + $r).Value = 'Hortaliza'
    }
}
